$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 15.9562
$ws.Range("E14").Value = 16.9293
$ws.Range("E21").Value = 16.7429
$ws.Range("E23").Value = 16.16799999999998
$ws.Range("E25").Value = 16.95920000000001
$ws.Range("E26").Value = 16.1946
$ws.Range("E29").Value = 17.02000000000002
$ws.Range("E53").Value = 16.62510000000002
$ws.Range("E57").Value = 16.7501
$ws.Range("E59").Value = 16.09729999999999
$ws.Range("E69").Value = 17.25990000000003
$ws.Range("E79").Value = 18.03910000000002
$ws.Range("E83").Value = 16.5116
$ws.Range("E91").Value = 18.38330000000002
$ws.Range("E93").Value = 17.42400000000002
